# Flight_Mission_Cycle.xlsx edit
# "Now reads previous flight mission cycle and adds to new flight mission cycle"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Flight Mission Cycle
#   Setting/No. of cycles table now reads: Light switch / Typing / Read_output
# ---------------------------------------------------------------
$wsFMC = $wb.Worksheets.Item("Flight Mission Cycle")
$wsFMC.Select()

$wsFMC.Range("A2").Value = "Light switch"
$wsFMC.Range("B2").Value = 2

$wsFMC.Range("A3").Value = "Typing"
$wsFMC.Range("B3").Value = 2

$wsFMC.Range("A4").Value = "Read_output"
$wsFMC.Range("B4").ClearContents()

# Column A needs to widen to fit the new "Read_output" setting label
$wsFMC.Columns("A").AutoFit()

$null = $wsFMC.Range("C7").Select()

# ---------------------------------------------------------------
# Sheet: Example - selection only
# ---------------------------------------------------------------
$wsExample = $wb.Worksheets.Item("Example")
$wsExample.Select()
$null = $wsExample.Range("C4").Select()

# ---------------------------------------------------------------
# Sheet: Typing
#   Type switches from triangle to sinosoidal, Max_RoM/Min_RoM updated
# ---------------------------------------------------------------
$wsTyping = $wb.Worksheets.Item("Typing")
$wsTyping.Select()

$wsTyping.Range("B4").Value = "sinosoidal"
$wsTyping.Range("C5").Value = 10
$wsTyping.Range("C6").Value = -10

$null = $wsTyping.Range("F8").Select()

# ---------------------------------------------------------------
# Sheet: Light switch
#   Min_RoM updated
# ---------------------------------------------------------------
$wsLightSwitch = $wb.Worksheets.Item("Light switch")
$wsLightSwitch.Select()

$wsLightSwitch.Range("C6").Value = -30

$null = $wsLightSwitch.Range("K13:L13").Select()

# ---------------------------------------------------------------
# Re-select Flight Mission Cycle as the active sheet/tab
# ---------------------------------------------------------------
$wsFMC.Select()
$null = $wsFMC.Range("C7").Select()
